$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.400.17"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.872.76"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7168"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.98"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07811"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3071"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.25"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08256"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").Value = "1.860.41"
$ws.Range("E12").Value = "  -1.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7235"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.243"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.70"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "29.477.12"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.855"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007877"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.34"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "2.133.92"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.780"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1549"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.49"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.005"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.34"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.932"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("E30").Value = "  -5.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.482"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.328"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.088"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05247"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.198"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7172"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.676"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01866"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.721"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "1.181.40"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9090"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.017"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.98"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4306"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.47"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5364"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.150"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.022"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.72%  "
